$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the raw input values; dependent formula cells (H5, H6, B23, D23, F23, H23)
# will recalculate automatically.
$ws.Range("D5").Value = 16
$ws.Range("F5").Value = 1.5
$ws.Range("B6").Value = 9

# Update the selected cell/active cell shown in the sheet view.
$ws.Range("B6").Select()
